$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(51, 1).Value = 0.310122
$ws.Cells.Item(51, 2).Value = 0.941775
$ws.Cells.Item(51, 3).Value = 0.961183
$ws.Cells.Item(51, 4).Value = 0.85755
$ws.Cells.Item(51, 5).Value = 2076.045710210585
$ws.Cells.Item(51, 6).Value = "query"
